# Energy.xlsx update - "updated AMGN and added others"
$wb = $excel.ActiveWorkbook
$wsEnergy = $wb.Worksheets.Item("Energy")
$wsPFIE = $wb.Worksheets.Item("PFIE")

# ---------------------------------------------------------------------------
# 1. Build out the "Energy" comparison table
#    NOTE: the order in which new text values are written below matters -
#    it controls the order new entries are appended to the shared string
#    table, so it mirrors the order the original author typed them in.
# ---------------------------------------------------------------------------

# Row 2: merged group headers
$wsEnergy.Range("Q2:V2").Merge()
$wsEnergy.Range("W2:AB2").Merge()
$wsEnergy.Range("Q2").Value = "EV/EPS"
$wsEnergy.Range("W2").Value = "EPS"

# Row 3: column headers
$wsEnergy.Range("C3").Value = "Name"
$wsEnergy.Range("D3").Value = "Ticker"
$wsEnergy.Range("E3").Value = "Price"
$wsEnergy.Range("F3").Value = "MC"
$wsEnergy.Range("G3").Value = "NC"
$wsEnergy.Range("H3").Value = "EV"
$wsEnergy.Range("I3").Value = "Update"
$wsEnergy.Range("J3").Value = "SO"
$wsEnergy.Range("K3").Value = "Last"
$wsEnergy.Range("L3").Value = "NPV"
$wsEnergy.Range("M3").Value = "Upside"
$wsEnergy.Range("N3").Value = "ROIC"
$wsEnergy.Range("O3").Value = "Terminal"
$wsEnergy.Range("P3").Value = "Discount"
$wsEnergy.Range("AC3").Value = "Founded"
$wsEnergy.Range("AD3").Value = "Location"

# Year headers (row 3) for the EV/EPS (Q:V) and EPS (W:AB) blocks
$wsEnergy.Range("Q3").Value = 2022
$wsEnergy.Range("R3").Value = 2023
$wsEnergy.Range("S3").Value = 2024
$wsEnergy.Range("T3").Value = 2025
$wsEnergy.Range("U3").Value = 2026
$wsEnergy.Range("V3").Value = 2027
$wsEnergy.Range("W3").Value = 2022
$wsEnergy.Range("X3").Value = 2023
$wsEnergy.Range("Y3").Value = 2024
$wsEnergy.Range("Z3").Value = 2025
$wsEnergy.Range("AA3").Value = 2026
$wsEnergy.Range("AB3").Value = 2027

# Companies already tracked - rows 6 through 9
# Row 6: New Fortress Energy
$wsEnergy.Range("C6").Value = "New Fortress Energy"
$wsEnergy.Range("D6").Value = "NFE"
# Row 7: Quanta Services
$wsEnergy.Range("C7").Value = "Quanta Services"
$wsEnergy.Range("D7").Value = "PWR"
# Row 8: PrimeEnergy Resources
$wsEnergy.Range("C8").Value = "PrimeEnergy Resources"
$wsEnergy.Range("D8").Value = "PNRG"
# Row 9: Talen Energy
$wsEnergy.Range("C9").Value = "Talen Energy"
$wsEnergy.Range("D9").Value = "TLN"

# Newly-added companies - rows 4 and 5
# Row 4: Exxon Mobile
$wsEnergy.Range("C4").Value = "Exxon Mobile"
$wsEnergy.Range("D4").Value = "XOM"
# Row 5: Chevron
$wsEnergy.Range("C5").Value = "Chevron"
$wsEnergy.Range("D5").Value = "CVX"

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------

# Column widths
$wsEnergy.Columns.Item(1).ColumnWidth = 2.6328125
$wsEnergy.Columns.Item(2).ColumnWidth = 2
$wsEnergy.Columns.Item(3).ColumnWidth = 21.36328125

# "Header" font (plain, non-bold Arial 10) applied to name/ticker/founded/location cells
$hdrRanges = @("C3:D3","AC3:AD3","C4:D4","C5:D5")
foreach ($addr in $hdrRanges) {
    $r = $wsEnergy.Range($addr)
    $r.Font.ThemeColor = 1
}

# Right-aligned numeric-label columns E:P across rows 1-5
$rightRanges = @("E1:P1","E2:J2","L2:P2","E3:P3","E4:P4","E5:P5")
foreach ($addr in $rightRanges) {
    $r = $wsEnergy.Range($addr)
    $r.Font.ThemeColor = 1
    $r.HorizontalAlignment = -4152
}

# Centre-aligned year columns Q:AB for rows 1, 3, 4, 5
$centerRanges = @("Q1:AB1","Q3:AB3","Q4:AB4","Q5:AB5")
foreach ($addr in $centerRanges) {
    $r = $wsEnergy.Range($addr)
    $r.Font.ThemeColor = 1
    $r.HorizontalAlignment = -4108
}

# Merged row-2 group headers also centred
$r = $wsEnergy.Range("Q2:AB2")
$r.Font.ThemeColor = 1
$r.HorizontalAlignment = -4108

# K2: placeholder date cell (right aligned, date number format)
$r = $wsEnergy.Range("K2")
$r.Font.ThemeColor = 1
$r.HorizontalAlignment = -4152
$r.NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# Freeze panes / view setup for Energy sheet
# ---------------------------------------------------------------------------
$wsEnergy.Activate()
$wsEnergy.Range("D4").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.Zoom = 145
$wsEnergy.Range("A6").Select()

# ---------------------------------------------------------------------------
# 2. PFIE sheet - Energy is now the active tab, restore PFIE's own selection
# ---------------------------------------------------------------------------
$wsPFIE.Range("L10").Select()
